$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header / request info fields
$ws.Range("B9").Value = "FECHA SOLICITUD:  2021-07-02"
$ws.Range("F9").Value = "AREA: Compras"
$ws.Range("B11").Value = "NOMBRE CENTRO DE COSTO: Centro de Diseño tecnologico industrial"
$ws.Range("B13").Value = "NOMBRE DE JEFE DE OFICINA O COORDINADOR DE AREA: jair"
$ws.Range("F13").Value = "CEDULA: 11111111"
$ws.Range("B15").Value = "NOMBRE DE SERVIDOR PÚBLICO A QUIEN SE LE ASIGNARA EL BIEN: jair"
$ws.Range("F15").Value = "CEDULA: 1111111"
$ws.Range("B17").Value = "CÓDIGO DE GRUPO O FICHA DE CARACTERIZACIÓN: 123456"

# Table row 21
$ws.Range("C21").Value = "Reprograf"
$ws.Range("D21").Value = "Cm - Centimetro"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "n/aa"
